# Sort_keys.xlsx edit: relabel ASSY-prefixed headers with a hyphen,
# add a data summary pass over column D (shared formulas), and
# update the view/selection state left behind by the editing session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Relabel the header row (row 1): insert a hyphen after "ASSY" ------
$ws.Range("B1").Value = "ASSY-1"
$ws.Range("C1").Value = "ASSY-37"
$ws.Range("D1").Value = "ASSY-77"
$ws.Range("E1").Value = "ASSY-79"
$ws.Range("F1").Value = "ASSY-116"
$ws.Range("G1").Value = "ASSY-156"
$ws.Range("H1").Value = "ASSY-158A"
$ws.Range("I1").Value = "ASSY-158B"
$ws.Range("J1").Value = "ASSY-196"
$ws.Range("K1").Value = "ASSY-236"
$ws.Range("L1").Value = "ASSY-276"
$ws.Range("M1").Value = "ASSY-316"
# A1 ("NanoZ") is left untouched.

# --- Data summary: recompute column D (the "=A<row>" mirror column) as a
# pair of shared formulas, same as Excel does when the column is re-entered
# and filled down in one continuous pass.
$ws.Range("D2:D33").Formula = "=A2"
$ws.Range("D34:D65").Formula = "=A34"

# --- Restore the view/selection state captured in the saved file --------
$ws.Range("M1").Select() | Out-Null
$wb.Windows.Item(1).Left = 1160
